$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Receipts table: update the Contributions value for receipt rec00001 ---
# (the way an expense's contributions are now recorded changed from a flat
# "paid" split, to tracking who actually paid whom)
$ws.Range("G8").Value = "Dan,£1,Dan/Marie,£1,Dan"

# --- Accounts table: recalculate/update the BalancesString for each account ---
# now that balances are correctly impacted by each new expense
$ws.Range("F3").Value = "Dan,3.00/Marie,-3.00"
$ws.Range("F4").Value = "Dan,1.00/Marie,2.00/Teddylou,-3.00"

# Widen column G so the longer Contributions/BalancesString values fit
$ws.Columns.Item(7).ColumnWidth = 21.5

# Update the selected cell / scroll position left over from editing
$ws.Range("F5").Select()
